$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values re-pulled from source data
$updates = @{
    2  = 2
    3  = 0
    4  = -2
    6  = 2
    7  = -1
    8  = -1
    9  = -1
    10 = -1
    11 = -3
    12 = -5
    13 = -4
    14 = 3
    15 = 4
    16 = -3
    17 = 5
    21 = -4
    22 = -3
    23 = -1
    24 = -2
    25 = 3
    26 = 3
    27 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
